# Updates the "phenotypic-feature-assertion-vs" ValueSet workbook:
#  1. Bumps the "Date" metadata value.
#  2. Rewrites the "Include #0" sheet from a CodeSystem-based "Codes / All
#     codes" listing to an explicit LOINC concept list (Concept/Description
#     rows for LA9633-4/Present and LA9634-2/Absent) plus the System URI.
#  3. Adds a new "Include ValueSet #1" sheet referencing the HL7 v3-NullFlavor
#     ValueSet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata!B8 (Date) bump
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-22T19:44:39+00:00"

$inc0 = $wb.Worksheets.Item("Include #0")

# ---------------------------------------------------------------------------
# 2. Add the "Include ValueSet #1" sheet -- duplicate "Include #0" *first*,
#    while its column B is still completely empty (only A1..A4 are
#    populated), so the new sheet inherits the exact column widths/styles
#    without dragging along any column-B cells. Then trim it down and rename
#    it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc0.Copy($null, $lastSheet)
$valueSetSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$valueSetSheet.Name = "Include ValueSet #1"

$valueSetSheet.Rows("3:4").Delete()

$valueSetSheet.Range("A1").Value = "ValueSet URL"
$valueSetSheet.Range("A2").Value = "http://terminology.hl7.org/ValueSet/v3-NullFlavor"

# ---------------------------------------------------------------------------
# 3. Rebuild "Include #0" sheet content
# ---------------------------------------------------------------------------

# Row 1 and row 2 previously only had column A populated -- B1/B2 are
# brand-new cells, so copy the existing bordered/header formatting onto them
# (from A1, and from the already-bordered A3) before writing the text.
$inc0.Range("A1").Copy()
$inc0.Range("B1").PasteSpecial(-4122)
$inc0.Range("A3").Copy()
$inc0.Range("B2").PasteSpecial(-4122)

$inc0.Range("A1").Value = "Concept"
$inc0.Range("B1").Value = "Description"

$inc0.Range("A2").Value = "LA9633-4"
$inc0.Range("B2").Value = "Present"

$inc0.Range("A3").Value = "LA9634-2"
$inc0.Range("B3").Value = "Absent"

$inc0.Range("A4").Value = ""
$inc0.Range("B4").Value = ""

# New row 5 -- copy the formatting from row 4 (already styled/bordered) then
# set the System URI values.
$inc0.Range("A4:B4").Copy()
$inc0.Range("A5:B5").PasteSpecial(-4122)
$inc0.Range("A5").Value = "System URI"
$inc0.Range("B5").Value = "http://loinc.org"

# Restore the original active sheet/tab selection.
$meta.Select()
